$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, shifting existing rows 192..262 down to 193..263
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new weekly price record
$ws.Range("A192").Value = 3
$ws.Range("B192").Value = 'Femacal de La Calera'
$ws.Range("C192").Value = 'Coquimbo'
$ws.Range("D192").Value = 44985
$ws.Range("E192").Value = 5
$ws.Range("F192").Value = 100112030
$ws.Range("G192").Value = 'Poroto granado'
$ws.Range("H192").Value = 'Sin especificar'
$ws.Range("I192").Value = 'Primera'
$ws.Range("J192").Value = 40
$ws.Range("K192").Value = 35000
$ws.Range("L192").Value = 35000
$ws.Range("M192").Value = 35000
$ws.Range("N192").Value = '$/malla 25 kilos'
$ws.Range("O192").Value = 'Provincia de Quillota'
$ws.Range("P192").Value = 1400
$ws.Range("Q192").Value = 25
$ws.Range("R192").Value = 'Hortaliza'
